$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Mes" date column (C) for all data rows from Aug 1, 2023 to Sep 1, 2023
$ws.Range("C2:C41").Value = 45170

# Update "Capacidad instalada" column (D) values that changed
$dUpdates = @{
    2 = 564
    3 = 312
    10 = 70
    11 = 240
    14 = 360
    15 = 424
    16 = 154
    17 = 1742
    18 = 161
    19 = 372
    20 = 372
    21 = 312
    23 = 86
    24 = 564
    27 = 392
    28 = 432
    30 = 372
    32 = 34
    33 = 550
}

foreach ($row in $dUpdates.Keys) {
    $ws.Cells.Item($row, 4).Value = $dUpdates[$row]
}

